$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$n = 20
$pMin = 0.000000001
$pMax = 100000

for ($i = 0; $i -lt $n; $i++) {
    $row = $i + 2
    $pressure = $pMin + ($pMax - $pMin) * $i / ($n - 1)
    $density = [Math]::Sqrt($pressure / 150)
    $ws.Cells.Item($row, 1).Value = $density
    $ws.Cells.Item($row, 2).Value = $pressure
}
